# Fix Bug in CVRP_L  & Add show_function
# Corrects the "class" (column E) values on Sheet1 that were computed
# with an off-by-something bug in the CVRP_L data generator.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$classFixes = @{
    4 = 10
    5 = 10
    6 = 10
    7 = 3
    12 = 3
    13 = 3
    14 = 3
    15 = 4
    16 = 4
    18 = 4
    19 = 4
    20 = 4
    22 = 4
    23 = 4
    24 = 4
    25 = 4
    26 = 4
    27 = 12
    28 = 12
    30 = 12
    32 = 12
    33 = 12
    34 = 12
    36 = 12
    38 = 12
    41 = 2
    42 = 2
    43 = 2
    50 = 2
    51 = 2
    53 = 2
    54 = 2
    55 = 8
    56 = 8
    57 = 8
    60 = 8
    62 = 8
    66 = 8
    69 = 13
    70 = 13
    72 = 13
    73 = 13
    74 = 13
    75 = -1
    76 = 10
    79 = 10
    80 = 10
    81 = 3
    82 = 10
    84 = 1
    85 = 1
    86 = 1
    87 = 1
    88 = 1
    89 = 1
    91 = 1
    92 = 1
    97 = 9
    98 = 9
    100 = 9
    101 = 9
    102 = 9
    103 = 9
    104 = 9
    105 = 9
    111 = -1
    114 = 11
    115 = 11
    116 = 8
    118 = 11
    119 = 0
    121 = 0
    122 = 0
    123 = 0
    126 = 0
    128 = 0
    129 = 0
    130 = 7
    131 = 7
    132 = 7
    133 = 7
    134 = 7
    135 = 7
    136 = 7
    137 = 7
    138 = 7
    139 = 7
    140 = 7
    141 = 7
    142 = 5
    147 = 5
    149 = 5
    151 = 5
    153 = 13
    154 = 13
    155 = 6
    156 = 6
    158 = 13
    159 = -1
    160 = -1
    161 = 13
    162 = 3
    163 = 3
    171 = 3
    172 = 10
}

foreach ($row in $classFixes.Keys) {
    $ws.Cells.Item($row, 5).Value = $classFixes[$row]
}
